$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 200
$ws.Range("J33").Value = 202
$ws.Range("L33").Value = 202
$ws.Range("N33").Value = -660

$ws.Range("H74").Value = 11592
$ws.Range("I74").Value = 3910.4
$ws.Range("J74").Value = 50000
$ws.Range("K74").Value = 3910.4
$ws.Range("L74").Value = 50000
$ws.Range("M74").Value = -2974.4
$ws.Range("N74").Value = -51872

$ws.Range("H77").Value = 11592
$ws.Range("I77").Value = 3910.4
$ws.Range("J77").Value = 50000
$ws.Range("K77").Value = 19552
$ws.Range("L77").Value = 250000
$ws.Range("M77").Value = -14872
$ws.Range("N77").Value = -259360

$ws.Range("H125").Value = 1832.6666
$ws.Range("I125").Value = 1998
$ws.Range("J125").Value = 1750
$ws.Range("K125").Value = 17982
$ws.Range("L125").Value = 15750
$ws.Range("M125").Value = -15522
$ws.Range("N125").Value = -20670

$ws.Range("H129").Value = 1263.7142
$ws.Range("I129").Value = 1433
$ws.Range("J129").Value = 1136.75
$ws.Range("K129").Value = 4299
$ws.Range("L129").Value = 3410.25
$ws.Range("M129").Value = 701
$ws.Range("N129").Value = -13410.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1207.3334
$ws.Range("I2").Value = 1135.2727
$ws.Range("K2").Value = 1135.2727
$ws.Range("M2").Value = -1022.2727

$ws.Range("H32").Value = 1145.2727
$ws.Range("I32").Value = 1059.8
$ws.Range("K32").Value = 1059.8
$ws.Range("M32").Value = -772.8

$ws.Range("H88").Value = 1321.8889
$ws.Range("I88").Value = 1326.8572
$ws.Range("J88").Value = 1304.5
$ws.Range("K88").Value = 1326.8572
$ws.Range("L88").Value = 1304.5
$ws.Range("M88").Value = -920.8571999999999
$ws.Range("N88").Value = -2116.5

$ws.Range("H91").Value = 1321.8889
$ws.Range("I91").Value = 1326.8572
$ws.Range("J91").Value = 1304.5
$ws.Range("K91").Value = 1326.8572
$ws.Range("L91").Value = 1304.5
$ws.Range("M91").Value = 77.14280000000008
$ws.Range("N91").Value = -4112.5

$ws.Range("H102").Value = 6314.3335
$ws.Range("I102").Value = 4471.5
$ws.Range("K102").Value = 4471.5
$ws.Range("M102").Value = -2849.5

$ws.Range("H110").Value = 1498.8235
$ws.Range("I110").Value = 1413.2307
$ws.Range("K110").Value = 1413.2307
$ws.Range("M110").Value = 631.7692999999999

$ws.Range("H116").Value = 1207.3334
$ws.Range("I116").Value = 1135.2727
$ws.Range("K116").Value = 1135.2727
$ws.Range("M116").Value = 1158.7273

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1207.3334
$ws.Range("I3").Value = 1135.2727
$ws.Range("K3").Value = 1135.2727
$ws.Range("M3").Value = -1021.2727

$ws.Range("H87").Value = 50000
$ws.Range("I87").Value = 50000
$ws.Range("K87").Value = 50000
$ws.Range("M87").Value = -48752

$ws.Range("H90").Value = 50000
$ws.Range("I90").Value = 50000
$ws.Range("K90").Value = 150000
$ws.Range("M90").Value = -143760

$ws.Range("H94").Value = 686.5714
$ws.Range("I94").Value = 686.5714
$ws.Range("K94").Value = 686.5714
$ws.Range("M94").Value = -235.5714

$ws.Range("H135").Value = 175516.67
$ws.Range("J135").Value = 175516.67
$ws.Range("L135").Value = 175516.67
$ws.Range("N135").Value = -185656.67

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H47").Value = 24900
$ws.Range("I47").Value = 24900
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 24900
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = -24334
$ws.Range("N47").ClearContents()

$ws.Range("H86").Value = 3701.3333
$ws.Range("I86").Value = 3299.5
$ws.Range("K86").Value = 3299.5
$ws.Range("M86").Value = -2176.5

$ws.Range("H89").Value = 3701.3333
$ws.Range("I89").Value = 3299.5
$ws.Range("K89").Value = 16497.5
$ws.Range("M89").Value = -10881.5

$ws.Range("H99").Value = 2597.6667
$ws.Range("J99").Value = 1394
$ws.Range("L99").Value = 1394
$ws.Range("N99").Value = -4390

$ws.Range("H126").Value = 2597.6667
$ws.Range("J126").Value = 1394
$ws.Range("L126").Value = 4182
$ws.Range("N126").Value = -9122

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 9835
$ws.Range("I80").Value = 9752.5
$ws.Range("K80").Value = 9752.5
$ws.Range("M80").Value = -8754.5

$ws.Range("H83").Value = 9835
$ws.Range("I83").Value = 9752.5
$ws.Range("K83").Value = 48762.5
$ws.Range("M83").Value = -43770.5

$ws.Range("H102").Value = 3102.75
$ws.Range("I102").Value = 2833.1428
$ws.Range("J102").Value = 4990
$ws.Range("K102").Value = 2833.1428
$ws.Range("L102").Value = 4990
$ws.Range("M102").Value = -1211.1428
$ws.Range("N102").Value = -8234

$ws.Range("H122").Value = 2831.3333
$ws.Range("I122").Value = 2831.3333
$ws.Range("K122").Value = 8493.999899999999
$ws.Range("M122").Value = -6043.999899999999

$ws.Range("H126").Value = 9900
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2225
$ws.Range("I22").Value = 2616.6667
$ws.Range("K22").Value = 2616.6667
$ws.Range("M22").Value = -2321.6667

$ws.Range("H27").Value = 2225
$ws.Range("I27").Value = 2616.6667
$ws.Range("K27").Value = 2616.6667
$ws.Range("M27").Value = -2509.6667

$ws.Range("H46").Value = 6727.273
$ws.Range("J46").Value = 6900
$ws.Range("L46").Value = 6900
$ws.Range("N46").Value = -7276

$ws.Range("H132").Value = 2373.5
$ws.Range("I132").Value = 2373.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7120.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4590.5
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 359
$ws.Range("I81").Value = 319
$ws.Range("K81").Value = 638
$ws.Range("M81").Value = 423

$ws.Range("H84").Value = 359
$ws.Range("I84").Value = 319
$ws.Range("K84").Value = 3190
$ws.Range("M84").Value = 2114

$ws.Range("H113").Value = 879.4375
$ws.Range("I113").Value = 792.8182
$ws.Range("K113").Value = 2378.4546
$ws.Range("M113").Value = -208.4546

